# Update cryptos list with latest prices / volume changes
# Commit: Updated cryptos list on Tue Aug 15 07:34:22 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "29.325.05"
$ws.Range("E2").Value = "  -0.22%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.838.94"
$ws.Range("E3").Value = "  -0.53%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.10%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.14"
$ws.Range("E5").Value = "  -0.56%  "

# Row 6 - XRP
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6254"
$ws.Range("E6").Value = "  -0.26%  "

# Row 7 - USDC
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  -0.05%  "

# Row 8 - Dogecoin
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07413"
$ws.Range("E8").Value = "  -1.12%  "

# Row 9 - Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2890"
$ws.Range("E9").Value = "  -0.37%  "

# Row 10 - Solana
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.73"
$ws.Range("E10").Value = "  +1.24%  "

# Row 11 - TRON
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07717"
$ws.Range("E11").Value = "  -0.32%  "

# Row 12 - WrappedEther
$ws.Range("D12").Value = "1.836.74"
$ws.Range("E12").Value = "  -0.64%  "

# Row 13 - Polkadot
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.946"
$ws.Range("E13").Value = "  -1.08%  "

# Row 14 - Polygon
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6728"
$ws.Range("E14").Value = "  -1.07%  "

# Row 15 - ShibaInu
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001017"
$ws.Range("E15").Value = "  -2.36%  "

# Row 16 - Litecoin
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "81.60"
$ws.Range("E16").Value = "  -0.67%  "

# Row 17 - Uniswap
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.205"
$ws.Range("E17").Value = "  +0.40%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "29.356.72"
$ws.Range("E18").Value = "  -0.26%  "

# Row 19 - BitcoinCash
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "231.71"
$ws.Range("E19").Value = "  +1.05%  "

# Row 20 - Avalanche
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.27"
$ws.Range("E20").Value = "  -0.60%  "

# Row 21 - Dai
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -0.03%  "

# Row 22 - Chainlink
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.326"
$ws.Range("E22").Value = "  -1.91%  "

# Row 23 - BinanceUSD
$ws.Range("E23").Value = "  +0.02%  "

# Row 24 - Monero
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "158.05"
$ws.Range("E24").Value = "  -0.38%  "

# Row 25 - Cosmos
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.456"
$ws.Range("E25").Value = "  +0.51%  "

# Row 26 - Stellar
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1344"
$ws.Range("E26").Value = "  -2.35%  "

# Row 27 - EthereumClassic
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.32"
$ws.Range("E27").Value = "  -1.28%  "

# Row 28 - Hedera
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.07252"
$ws.Range("E28").Value = "  +12.90%  "

# Row 29 - Toncoin
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.457"
$ws.Range("E29").Value = "  +4.88%  "

# Row 30 - PancakeSwap
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.478"
$ws.Range("E30").Value = "  +0.31%  "

# Row 31 - Filecoin
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.034"
$ws.Range("E31").Value = "  -1.51%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.034"
$ws.Range("E32").Value = "  -0.89%  "

# Row 33 - LidoDAOToken
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.814"
$ws.Range("E33").Value = "  -0.80%  "

# Row 35 - ImmutableX
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6941"
$ws.Range("E35").Value = "  -0.86%  "

# Row 36 - HuobiToken
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.569"
$ws.Range("E36").Value = "  -0.41%  "

# Row 37 - VeChain -> FraxShare
$ws.Range("B37").Value = "FraxShare"
$ws.Range("C37").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.925"
$ws.Range("E37").Value = "  +4.66%  "

# Row 38 - FraxShare -> VeChain
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01830"
$ws.Range("E38").Value = "  -0.06%  "

# Row 39 - MXToken
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.810"
$ws.Range("E39").Value = "  -0.90%  "

# Row 40 - Maker
$ws.Range("D40").Value = "1.228.24"
$ws.Range("E40").Value = "  -2.88%  "

# Row 41 - TrustWalletToken
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9426"
$ws.Range("E41").Value = "  +3.51%  "

# Row 42 - PaxDollar
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9998"
$ws.Range("E42").Value = "  -0.03%  "

# Row 43 - RocketPoolETH
$ws.Range("D43").Value = "1.990.58"
$ws.Range("E43").Value = "  -0.95%  "

# Row 44 - Quant
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.36"
$ws.Range("E44").Value = "  -1.21%  "

# Row 45 - Aave
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.31"
$ws.Range("E45").Value = "  -1.54%  "

# Row 46 - BabyDogeCoin
$ws.Range("E46").Value = "  +6.77%  "

# Row 47 - RenderToken
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.698"
$ws.Range("E47").Value = "  -3.39%  "

# Row 48 - Aptos
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.921"
$ws.Range("E48").Value = "  -2.40%  "

# Row 49 - Algorand -> EnergySwap
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.880"
$ws.Range("E49").Value = "  -1.88%  "

# Row 50 - EnergySwap -> Algorand
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1133"
$ws.Range("E50").Value = "  -3.44%  "

# Row 51 - TheSandbox
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3896"
$ws.Range("E51").Value = "  -1.41%  "
